$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44316
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 17500
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 17750
$ws.Range("Q3").Value = '$/caja 16 kilos granel'
$ws.Range("S3").Value = 1109
$ws.Range("T3").Value = 16

# Row 4
$ws.Range("D4").Value = 44316
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 16000
$ws.Range("Q4").Value = '$/caja 16 kilos granel'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 16

# Row 5
$ws.Range("D5").Value = 44334
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 12500
$ws.Range("Q5").Value = '$/caja 12 kilos empedrada'
$ws.Range("S5").Value = 1042
$ws.Range("T5").Value = 12

# Row 6
$ws.Range("D6").Value = 45085
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 19000
$ws.Range("P6").Value = 18500
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Región del Maule'
$ws.Range("S6").Value = 1028
$ws.Range("T6").Value = 18

# Row 8
$ws.Range("D8").NumberFormat = $ws.Range("D7").NumberFormat
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C8").Value = 'Ñuble'
$ws.Range("D8").Value = 44344
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 'Fruta'
$ws.Range("G8").Value = 100107
$ws.Range("H8").Value = 'Otros'
$ws.Range("I8").Value = 100107001
$ws.Range("J8").Value = 'Caqui'
$ws.Range("K8").Value = 'Mankaki'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 120
$ws.Range("N8").Value = 13000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 13500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Provincia de Curicó'
$ws.Range("S8").Value = 750
$ws.Range("T8").Value = 18
